# Update NATMI LR-pair sheet (Vegfa-Nrp1) with recomputed TPM-based statistics.
# Ligand/receptor average & total expression, their derived specificities, and the
# edge-level weights/specificities (which are simple products of the ligand- and
# receptor-side values) are refreshed to match the new TPM input used by the
# upstream scripts ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.062550666666667
$ws.Range("H2").Value = 9.187652
$ws.Range("I2").Value = 0.06849600470812313
$ws.Range("J2").Value = 0.06849600470812313
$ws.Range("M2").Value = 86.89540866666668
$ws.Range("N2").Value = 260.686226
$ws.Range("O2").Value = 0.319779657009892
$ws.Range("P2").Value = 0.3197796570098919
$ws.Range("Q2").Value = 266.1215917423725
$ws.Range("R2").Value = 2395.094325681352
$ws.Range("S2").Value = 0.02190362889211156
$ws.Range("T2").Value = 0.02190362889211156
$ws.Range("G3").Value = 3.062550666666667
$ws.Range("H3").Value = 9.187652
$ws.Range("I3").Value = 0.06849600470812313
$ws.Range("J3").Value = 0.06849600470812313
$ws.Range("O3").Value = 0.1999969065479545
$ws.Range("P3").Value = 0.1999969065479545
$ws.Range("Q3").Value = 166.4380267705578
$ws.Range("R3").Value = 1497.94224093502
$ws.Range("S3").Value = 0.01369898905251876
$ws.Range("T3").Value = 0.01369898905251875
$ws.Range("G4").Value = 3.062550666666667
$ws.Range("H4").Value = 9.187652
$ws.Range("I4").Value = 0.06849600470812313
$ws.Range("J4").Value = 0.06849600470812313
$ws.Range("M4").Value = 60.92601633333334
$ws.Range("N4").Value = 182.778049
$ws.Range("O4").Value = 0.224210932487692
$ws.Range("P4").Value = 0.224210932487692
$ws.Range("Q4").Value = 186.5890119389942
$ws.Range("R4").Value = 1679.301107450948
$ws.Range("S4").Value = 0.01535755308728963
$ws.Range("T4").Value = 0.01535755308728963
$ws.Range("G5").Value = 3.062550666666667
$ws.Range("H5").Value = 9.187652
$ws.Range("I5").Value = 0.06849600470812313
$ws.Range("J5").Value = 0.06849600470812313
$ws.Range("M5").Value = 7.809668333333332
$ws.Range("N5").Value = 23.429005
$ws.Range("O5").Value = 0.02873998867505581
$ws.Range("P5").Value = 0.02873998867505581
$ws.Range("Q5").Value = 23.91750496069555
$ws.Range("R5").Value = 215.25754464626
$ws.Range("S5").Value = 0.001968574399598029
$ws.Range("T5").Value = 0.001968574399598029
$ws.Range("G6").Value = 3.062550666666667
$ws.Range("H6").Value = 9.187652
$ws.Range("I6").Value = 0.06849600470812313
$ws.Range("J6").Value = 0.06849600470812313
$ws.Range("M6").Value = 61.75795633333333
$ws.Range("N6").Value = 185.273869
$ws.Range("O6").Value = 0.2272725152794058
$ws.Range("P6").Value = 0.2272725152794058
$ws.Range("Q6").Value = 189.1368703406209
$ws.Range("R6").Value = 1702.231833065588
$ws.Range("S6").Value = 0.01556725927660517
$ws.Range("T6").Value = 0.01556725927660517
$ws.Range("I7").Value = 0.396815038797359
$ws.Range("J7").Value = 0.396815038797359
$ws.Range("M7").Value = 86.89540866666668
$ws.Range("N7").Value = 260.686226
$ws.Range("O7").Value = 0.319779657009892
$ws.Range("P7").Value = 0.3197796570098919
$ws.Range("Q7").Value = 1541.711085223938
$ws.Range("R7").Value = 13875.39976701544
$ws.Range("S7").Value = 0.1268933770029864
$ws.Range("T7").Value = 0.1268933770029864
$ws.Range("I8").Value = 0.396815038797359
$ws.Range("J8").Value = 0.396815038797359
$ws.Range("O8").Value = 0.1999969065479545
$ws.Range("P8").Value = 0.1999969065479545
$ws.Range("S8").Value = 0.07936178023117837
$ws.Range("T8").Value = 0.07936178023117836
$ws.Range("I9").Value = 0.396815038797359
$ws.Range("J9").Value = 0.396815038797359
$ws.Range("M9").Value = 60.92601633333334
$ws.Range("N9").Value = 182.778049
$ws.Range("O9").Value = 0.224210932487692
$ws.Range("P9").Value = 0.224210932487692
$ws.Range("Q9").Value = 1080.958317601729
$ws.Range("R9").Value = 9728.62485841556
$ws.Range("S9").Value = 0.08897026987389554
$ws.Range("T9").Value = 0.08897026987389553
$ws.Range("I10").Value = 0.396815038797359
$ws.Range("J10").Value = 0.396815038797359
$ws.Range("M10").Value = 7.809668333333332
$ws.Range("N10").Value = 23.429005
$ws.Range("O10").Value = 0.02873998867505581
$ws.Range("P10").Value = 0.02873998867505581
$ws.Range("Q10").Value = 138.5602809880222
$ws.Range("R10").Value = 1247.0425288922
$ws.Range("S10").Value = 0.01140445972112793
$ws.Range("T10").Value = 0.01140445972112793
$ws.Range("I11").Value = 0.396815038797359
$ws.Range("J11").Value = 0.396815038797359
$ws.Range("M11").Value = 61.75795633333333
$ws.Range("N11").Value = 185.273869
$ws.Range("O11").Value = 0.2272725152794058
$ws.Range("P11").Value = 0.2272725152794058
$ws.Range("Q11").Value = 1095.718719099596
$ws.Range("R11").Value = 9861.468471896358
$ws.Range("S11").Value = 0.0901851519681708
$ws.Range("T11").Value = 0.09018515196817078
$ws.Range("G12").Value = 13.27534766666667
$ws.Range("H12").Value = 39.826043
$ws.Range("I12").Value = 0.2969120759943797
$ws.Range("J12").Value = 0.2969120759943796
$ws.Range("M12").Value = 86.89540866666668
$ws.Range("N12").Value = 260.686226
$ws.Range("O12").Value = 0.319779657009892
$ws.Range("P12").Value = 0.3197796570098919
$ws.Range("Q12").Value = 1153.56676068708
$ws.Range("R12").Value = 10382.10084618372
$ws.Range("S12").Value = 0.09494644182357771
$ws.Range("T12").Value = 0.09494644182357768
$ws.Range("G13").Value = 13.27534766666667
$ws.Range("H13").Value = 39.826043
$ws.Range("I13").Value = 0.2969120759943797
$ws.Range("J13").Value = 0.2969120759943796
$ws.Range("O13").Value = 0.1999969065479545
$ws.Range("P13").Value = 0.1999969065479545
$ws.Range("Q13").Value = 721.4648542412561
$ws.Range("R13").Value = 6493.183688171304
$ws.Range("S13").Value = 0.05938149671560713
$ws.Range("T13").Value = 0.05938149671560711
$ws.Range("G14").Value = 13.27534766666667
$ws.Range("H14").Value = 39.826043
$ws.Range("I14").Value = 0.2969120759943797
$ws.Range("J14").Value = 0.2969120759943796
$ws.Range("M14").Value = 60.92601633333334
$ws.Range("N14").Value = 182.778049
$ws.Range("O14").Value = 0.224210932487692
$ws.Range("P14").Value = 0.224210932487692
$ws.Range("Q14").Value = 808.8140487700119
$ws.Range("R14").Value = 7279.326438930108
$ws.Range("S14").Value = 0.06657093342555634
$ws.Range("T14").Value = 0.06657093342555633
$ws.Range("G15").Value = 13.27534766666667
$ws.Range("H15").Value = 39.826043
$ws.Range("I15").Value = 0.2969120759943797
$ws.Range("J15").Value = 0.2969120759943796
$ws.Range("M15").Value = 7.809668333333332
$ws.Range("N15").Value = 23.429005
$ws.Range("O15").Value = 0.02873998867505581
$ws.Range("P15").Value = 0.02873998867505581
$ws.Range("Q15").Value = 103.6760622863572
$ws.Range("R15").Value = 933.0845605772148
$ws.Range("S15").Value = 0.008533249701565784
$ws.Range("T15").Value = 0.008533249701565783
$ws.Range("G16").Value = 13.27534766666667
$ws.Range("H16").Value = 39.826043
$ws.Range("I16").Value = 0.2969120759943797
$ws.Range("J16").Value = 0.2969120759943796
$ws.Range("M16").Value = 61.75795633333333
$ws.Range("N16").Value = 185.273869
$ws.Range("O16").Value = 0.2272725152794058
$ws.Range("P16").Value = 0.2272725152794058
$ws.Range("Q16").Value = 819.8583415078185
$ws.Range("R16").Value = 7378.725073570366
$ws.Range("S16").Value = 0.06747995432807277
$ws.Range("T16").Value = 0.06747995432807274
$ws.Range("G17").Value = 3.455866
$ws.Range("H17").Value = 10.367598
$ws.Range("I17").Value = 0.07729276657626213
$ws.Range("J17").Value = 0.07729276657626213
$ws.Range("M17").Value = 86.89540866666668
$ws.Range("N17").Value = 260.686226
$ws.Range("O17").Value = 0.319779657009892
$ws.Range("P17").Value = 0.3197796570098919
$ws.Range("Q17").Value = 300.2988883672387
$ws.Range("R17").Value = 2702.689995305148
$ws.Range("S17").Value = 0.02471665438510275
$ws.Range("T17").Value = 0.02471665438510274
$ws.Range("G18").Value = 3.455866
$ws.Range("H18").Value = 10.367598
$ws.Range("I18").Value = 0.07729276657626213
$ws.Range("J18").Value = 0.07729276657626213
$ws.Range("O18").Value = 0.1999969065479545
$ws.Range("P18").Value = 0.1999969065479545
$ws.Range("Q18").Value = 187.8132251276367
$ws.Range("R18").Value = 1690.31902614873
$ws.Range("S18").Value = 0.01545831421378556
$ws.Range("T18").Value = 0.01545831421378556
$ws.Range("G19").Value = 3.455866
$ws.Range("H19").Value = 10.367598
$ws.Range("I19").Value = 0.07729276657626213
$ws.Range("J19").Value = 0.07729276657626213
$ws.Range("M19").Value = 60.92601633333334
$ws.Range("N19").Value = 182.778049
$ws.Range("O19").Value = 0.224210932487692
$ws.Range("P19").Value = 0.224210932487692
$ws.Range("Q19").Value = 210.5521483618113
$ws.Range("R19").Value = 1894.969335256302
$ws.Range("S19").Value = 0.01732988326861725
$ws.Range("T19").Value = 0.01732988326861724
$ws.Range("G20").Value = 3.455866
$ws.Range("H20").Value = 10.367598
$ws.Range("I20").Value = 0.07729276657626213
$ws.Range("J20").Value = 0.07729276657626213
$ws.Range("M20").Value = 7.809668333333332
$ws.Range("N20").Value = 23.429005
$ws.Range("O20").Value = 0.02873998867505581
$ws.Range("P20").Value = 0.02873998867505581
$ws.Range("Q20").Value = 26.98916726444333
$ws.Range("R20").Value = 242.9025053799899
$ws.Range("S20").Value = 0.002221393236065506
$ws.Range("T20").Value = 0.002221393236065506
$ws.Range("G21").Value = 3.455866
$ws.Range("H21").Value = 10.367598
$ws.Range("I21").Value = 0.07729276657626213
$ws.Range("J21").Value = 0.07729276657626213
$ws.Range("M21").Value = 61.75795633333333
$ws.Range("N21").Value = 185.273869
$ws.Range("O21").Value = 0.2272725152794058
$ws.Range("P21").Value = 0.2272725152794058
$ws.Range("Q21").Value = 213.4272215218513
$ws.Range("R21").Value = 1920.844993696662
$ws.Range("S21").Value = 0.01756652147269109
$ws.Range("T21").Value = 0.01756652147269108
$ws.Range("G22").Value = 7.175465666666668
$ws.Range("H22").Value = 21.526397
$ws.Range("I22").Value = 0.1604841139238761
$ws.Range("J22").Value = 0.1604841139238761
$ws.Range("M22").Value = 86.89540866666668
$ws.Range("N22").Value = 260.686226
$ws.Range("O22").Value = 0.319779657009892
$ws.Range("P22").Value = 0.3197796570098919
$ws.Range("Q22").Value = 623.515021478636
$ws.Range("R22").Value = 5611.635193307723
$ws.Range("S22").Value = 0.05131955490611352
$ws.Range("T22").Value = 0.05131955490611351
$ws.Range("G23").Value = 7.175465666666668
$ws.Range("H23").Value = 21.526397
$ws.Range("I23").Value = 0.1604841139238761
$ws.Range("J23").Value = 0.1604841139238761
$ws.Range("O23").Value = 0.1999969065479545
$ws.Range("P23").Value = 0.1999969065479545
$ws.Range("Q23").Value = 389.9593759275662
$ws.Range("R23").Value = 3509.634383348096
$ws.Range("S23").Value = 0.03209632633486473
$ws.Range("T23").Value = 0.03209632633486473
$ws.Range("G24").Value = 7.175465666666668
$ws.Range("H24").Value = 21.526397
$ws.Range("I24").Value = 0.1604841139238761
$ws.Range("J24").Value = 0.1604841139238761
$ws.Range("M24").Value = 60.92601633333334
$ws.Range("N24").Value = 182.778049
$ws.Range("O24").Value = 0.224210932487692
$ws.Range("P24").Value = 0.224210932487692
$ws.Range("Q24").Value = 437.172538406606
$ws.Range("R24").Value = 3934.552845659454
$ws.Range("S24").Value = 0.03598229283233325
$ws.Range("T24").Value = 0.03598229283233324
$ws.Range("G25").Value = 7.175465666666668
$ws.Range("H25").Value = 21.526397
$ws.Range("I25").Value = 0.1604841139238761
$ws.Range("J25").Value = 0.1604841139238761
$ws.Range("M25").Value = 7.809668333333332
$ws.Range("N25").Value = 23.429005
$ws.Range("O25").Value = 0.02873998867505581
$ws.Range("P25").Value = 0.02873998867505581
$ws.Range("Q25").Value = 56.03800699388722
$ws.Range("R25").Value = 504.342062944985
$ws.Range("S25").Value = 0.004612311616698565
$ws.Range("T25").Value = 0.004612311616698565
$ws.Range("G26").Value = 7.175465666666668
$ws.Range("H26").Value = 21.526397
$ws.Range("I26").Value = 0.1604841139238761
$ws.Range("J26").Value = 0.1604841139238761
$ws.Range("M26").Value = 61.75795633333333
$ws.Range("N26").Value = 185.273869
$ws.Range("O26").Value = 0.2272725152794058
$ws.Range("P26").Value = 0.2272725152794058
$ws.Range("Q26").Value = 443.1420953133326
$ws.Range("R26").Value = 3988.278857819993
$ws.Range("S26").Value = 0.03647362823386603
$ws.Range("T26").Value = 0.03647362823386602
